$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor (name unchanged) - update metric values
$ws.Range("B3").Value = 0.9997908077299428
$ws.Range("C3").Value = 0.9998092568769338
$ws.Range("D3").Value = 0.9977913833509873

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9997928356178981
$ws.Range("C4").Value = 0.9998053906431603
$ws.Range("D4").Value = 0.9991552353872226

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9998167229049079
$ws.Range("C5").Value = 0.9998096953785893
$ws.Range("D5").Value = 0.999797275258466
